# "delete task before creating new ones" - refreshed timestamps/dates for
# the re-run of the small-task scheduler (rows 2, 8, 9, 10 of Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Check_time only
$ws.Range("C2").Value = "2024-08-23 17:23:09"

# Row 8
$ws.Range("B8").Value = "task_2024-08-25_ZBZ_SMALL_READY"
$ws.Range("C8").Value = "2024-08-23 17:25:57"
# F8 holds a plain date-like string ("2024-08-25"); Excel would otherwise
# auto-convert that to a date serial number on input, so force Text first
# and restore the default (unstyled) look afterwards.
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2024-08-25"
$ws.Range("F8").Style = "Normal"

# Row 9
$ws.Range("B9").Value = "task_2024-08-28_ZBZ3_SMALL_READY"
$ws.Range("C9").Value = "2024-08-23 17:27:03"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "2024-08-28"
$ws.Range("F9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = "task_2024-08-23_HSG_SMALL_DONE"
$ws.Range("C10").Value = "2024-08-23 17:27:32"
$ws.Range("D10").Value = "2024-08-23 17:27:40"
$ws.Range("E10").Value = "2024-08-23 17:27:52"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "2024-08-23"
$ws.Range("F10").Style = "Normal"
